# Add the newest Microsoft Forms response (Mark Zbaracki) as row 12 of the
# response table on Sheet1, and grow Table1 to include it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new response row -------------------------------------------------
$ws.Range("A12").Value = 14
$ws.Range("B12").Value = 44027.4429282407
$ws.Range("B12").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("C12").Value = 44027.4474652778
$ws.Range("C12").NumberFormat = "m/d/yy\ h:mm:ss"
$ws.Range("D12").Value = "mzbarack@uwo.ca"
$ws.Range("E12").Value = "Mark Zbaracki"
$ws.Range("F12").Value = "Petrella;Smye;Tang;Esses;Lee;De Groot;Bitsuamlak;McNair;Beveridge;Hill;"

# --- grow the table to cover the new row -------------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F12"))
